# RPA datasets push 2024-05-30
# Insert a new IPO-subscription row for "미래에셋비전스팩4호" at the top of the
# data block (row 2), pushing the existing 13 rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert a new row above the current row 2 (shifts rows 2..14 -> 3..15)
$ws.Rows.Item(2).Insert()
# The inserted row inherits formatting from its neighbour; the rest of the
# data rows in this sheet carry no explicit style, so strip it back off.
$ws.Range("A2:T2").ClearFormats()

# Populate the new row 2 with the new record.
# Dates are stored as plain text in this sheet (not Excel date serials), so
# force the date-shaped cells to Text format before assigning.
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("A2").Value = "2024-05-20"
$ws.Range("B2").Value = "미래에셋비전스팩4호"
$ws.Range("C2").Value = "미래"
$ws.Range("D2").Value = "2024-05-23"
$ws.Range("E2").Value = "2024-05-29"
$ws.Range("F2").Value = 13300000
$ws.Range("G2").Value = 6650000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "687.21 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"
